# Update the SEBI tracker workbook: "Informal Guidance" entry -> "Regulations" (IFSCA) entry
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet / sheet tab
$ws.Name = "Regulations"

# Row 2 field updates
$ws.Range("B2").Value = "Regulations"
$ws.Range("C2").Value = 2026
$ws.Range("D2").Value = "February"
# Prefix with an apostrophe so Excel stores the literal text "2026-02-02"
# instead of auto-converting it to a date serial value (matches the
# original cell, which is a plain inline string, not a date).
$ws.Range("E2").Value = "'2026-02-02"
$ws.Range("F2").Value = 'IFSCA (Fund Management) (Amendment) Regulations, 2026
'
$ws.Range("G2").Value = 'https://ifsca.gov.in/CommonDirect/DownloadFile?id=36ff47aaeb9222f627d166fe86841979&fileName=106_IFSCA__Fund_Management___Amendment__Regulations__2026_20260202_0720.pdf'
$ws.Range("H2").Value = '106_IFSCA__Fund_Management___Amendment__Regulations__2026_20260202_0720.pdf'
$ws.Range("I2").Value = '/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Regulations/2026/February/106_IFSCA__Fund_Management___Amendment__Regulations__2026_20260202_0720.pdf'
$ws.Range("J2").Value = '• SEBI has amended the International Financial Services Centres Authority (Fund Management) Regulations, 2025 with the addition of the International Financial Services Centres Authority (Fund Management) (Amendment) Regulations, 2026.
• Effective from publication in the Official Gazette, these regulations introduce changes to improve governance and investor protection.
• Key updates include:
   - Increased experience requirements for Key Management Personnel (KMP): At least five years for all KMPs, except those referred under sub-regulation (2), who now require a minimum of three years of experience if they possess a professional qualification. Individuals with post-qualification experience of at least two or three years in an eligible institution may also be eligible for certain KMP roles.
   - Minimum corpus requirements: Fund Management Entities (FMEs) must now achieve a minimum size of corpus of USD 3 Million before investing in unlisted securities, and if they fail to do so within the specified time period, they may extend the validity of their placement memorandum for six months at a time, with each extension incurring increasing fees.
   - Appointment of custodian: FMEs that are required to appoint a custodian must now make this appointment within 24 months from the date of the new regulations'' commencement. During this period, they may appoint an independent custodian regulated by a foreign financial sector regulator.
• These changes aim to ensure greater accountability, transparency, and enhanced investor protection in the fund management industry under the International Financial Services Centres Authority (IFSCA).'
$ws.Range("K2").Value = '(1) These regulations may be called the International Financial Services Centres Authority (Fund Management) (Amendment) Regulations, 2026.
(2) These regulations shall come into force on the date of their publication in the Official Gazette.
In regulation 7 of the principal regulations, in sub-regulation (5), for clause (b), the following clause shall be substituted, namely: -
“In addition to the qualifications mentioned under clause (a), an experience of at least five (5) years in related activities in the securities market or financial products in an eligible institution:
Provided that for the KMP referred under sub-regulation (2), the experience mentioned in clause (b) shall be required for a minimum period of 3 (three) years, if such KMP possesses a professional qualification:
Provided further that an individual with post-qualification experience of at least 2 (two) years in an eligible institution in IFSC, India or any foreign jurisdiction and who holds valid certification(s), as specified by the Authority, shall be considered eligible for the KMP referred under sub-regulation (2):
Provided also that individuals with a post-qualification experience of at least 3 (three) years in an eligible institution in IFSC, India or any foreign jurisdiction and who holds valid certification(s), as specified by the Authority, shall be considered eligible for the KMP referred under sub-regulations (1), (3) and (4).
Explanation. – For the purpose of this clause, “eligible institution” shall include the following -
i) Market Infrastructure Institutions, Capital Market Intermediaries, financial sector regulators, FMEs, Banks, Finance Companies, Insurance Companies, and Insurance Intermediaries in IFSC, and equivalent institutions in India or any foreign jurisdiction;
ii) consulting firms / advisory firms / firms of Chartered Accountants / Company Secretaries / Cost Accountants in IFSC, India or any foreign jurisdiction, providing services to the institutions mentioned above in (i), in relation to a financial product; and
iii) a company, whether private or public, if the experience is in relation to finance/ accounts/ secretarial/ law departments of such company.
In regulation 19 of the principal regulations, in sub-regulation (3), for the proviso, the following proviso shall be substituted, namely: -
“Provided that if a FME fails to achieve the minimum size of corpus, as specified under sub-regulation (1) of regulation 23, within the specified time period, it shall have the option to extend the validity of the placement memorandum, wherein each such extension shall be for a period of six (6) months starting from the day after the expiry of the existing validity of the placement memorandum, by filing an application at such time when the placement memorandum is still valid, accompanied by a fee equal to -
i) for the first extension, twenty-five per cent. (25%) of the applicable fee for filing of a fresh scheme, as may be prevalent at the time of such extension; and
ii) for each subsequent extension, fifty per cent. (50%) of the applicable fee for filing of a fresh scheme, as may be prevalent at the time of such extension.”
In regulation 35 of the principal regulations,
in sub-regulation (1), after the proviso, the following proviso shall be inserted, namely: -
[भाग III—खण् ड 4]
भारत का रािपत्र : असाधारण
“Provided further that the investments by an open-ended scheme in unlisted securities shall be undertaken only upon achieving the minimum corpus of USD 3 Million.”
n sub-regulation (2), after the proviso, the following proviso shall be inserted, namely: -
“Provided further that if a FME fails to achieve the minimum size of corpus of USD 3 Million within the specified time period, it shall have the option to extend the validity of the placement memorandum, wherein each such extension shall be for a period of six (6) months starting from the day after the expiry of the existing validity of the placement memorandum, by filing an application at such time when the placement memorandum is still valid, accompanied by a fee equal to -
i) for the first extension, twenty-five per cent. (25%) of the applicable fee for filing of a fresh scheme, as may be prevalent at the time of such extension; and
ii) for each subsequent extension, fifty per cent. (50%) of the applicable fee for filing of a fresh scheme, as may be prevalent at the time of such extension”
In regulation 131 of the principal regulations, in sub-regulation (1), after clause (b), the following clauses shall be inserted, namely: -
“(c) If the FME has raised funds from the investors under the scheme but fails to achieve the minimum corpus during the validity or extended validity of the placement memorandum or offer document, as applicable, and the FME has not extended the validity thereof by making the requisite filing and payment of fee to the Authority.
In regulation 132 of the principal regulations, for the Explanation II, the following Explanation shall be substituted, namely: -
“Explanation II. – In case of schemes which are required to appoint custodian in IFSC in terms of the abovementioned provision, such appointment may be made within twenty four (24) months from the date of commencement of the International Financial Services Centres Authority (Fund Management) (Amendment) Regulation, 2026, during which period the FME may appoint an independent custodian in India, or any foreign jurisdiction, which is regulated by the financial sector regulator in that jurisdiction and make necessary arrangement to provide such information to Authority whenever directed to do so.”
PRAVEEN TRIVEDI, Executive Director
[ADVT.-III/4/Exty./640/2025-26]
1. The International Financial Services Centres Authority (Fund Management) Regulations, 2025, the principal regulations, were published in the Gazette of India on February 13, 2025, vide F. No. IFSCA/GN/2025/002.
International Financial Services Centres Authority (Fund Management) (Amendment) Regulations, 2025, were published in the Gazette of India on July 30, 2025, vide F. No. IFSCA/GN/2025/007.
Uploaded by Dte. of Printing at Government of India Press, Ring Road, Mayapuri, New Delhi-110064 and Published by the Controller of Publications, Delhi-110054.'
